# Slit1-Gpc1.xlsx: refresh the LR-pair table with new TPM-derived numbers.
#
# The "ECs" sending-cluster block (old rows 2-4) is dropped entirely, and the
# remaining "FAPs" / "MuSCs" sending-cluster blocks (old rows 5-10) move up
# to become the new rows 2-7, carrying updated statistics. Net effect: the
# sheet shrinks from 9 data rows (A1:T10) to 6 data rows (A1:T7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three trailing rows so the used range collapses back to A1:T7
# once the remaining rows below are rewritten in place.
$ws.Rows("8:10").Delete()

# New content for data rows 2-7 (columns A-T), keyed by row number.
$rows = @{
  2 = @{
    "A" = "FAPs"; "B" = "Slit1"; "C" = "Gpc1"; "D" = "ECs"
    "E" = 2; "F" = 0.6666666666666666; "G" = 0.2055996666666667; "H" = 0.616799
    "I" = 0.9059768423248155; "J" = 0.9059768423248156
    "K" = 3; "L" = 1
    "M" = 0.1112926666666667; "N" = 0.333878
    "O" = 0.01397697460904174; "P" = 0.01397697460904174
    "Q" = 0.02288173516911111; "R" = 0.205935616522
    "S" = 0.01266281532155375; "T" = 0.01266281532155375
  }
  3 = @{
    "A" = "FAPs"; "B" = "Slit1"; "C" = "Gpc1"; "D" = "FAPs"
    "E" = 2; "F" = 0.6666666666666666; "G" = 0.2055996666666667; "H" = 0.616799
    "I" = 0.9059768423248155; "J" = 0.9059768423248156
    "K" = 3; "L" = 1
    "M" = 3.316850333333333; "N" = 9.950551000000001
    "O" = 0.4165551449121381; "P" = 0.4165551449121381
    "Q" = 0.6819433229165556; "R" = 6.137489906249001
    "S" = 0.3773893148416548; "T" = 0.3773893148416549
  }
  4 = @{
    "A" = "FAPs"; "B" = "Slit1"; "C" = "Gpc1"; "D" = "MuSCs"
    "E" = 2; "F" = 0.6666666666666666; "G" = 0.2055996666666667; "H" = 0.616799
    "I" = 0.9059768423248155; "J" = 0.9059768423248156
    "K" = 3; "L" = 1
    "M" = 4.534429; "N" = 13.603287
    "O" = 0.5694678804788202; "P" = 0.5694678804788201
    "Q" = 0.9322770909236666; "R" = 8.390493818312999
    "S" = 0.515924712161607; "T" = 0.515924712161607
  }
  5 = @{
    "A" = "MuSCs"; "B" = "Slit1"; "C" = "Gpc1"; "D" = "ECs"
    "E" = 2; "F" = 0.6666666666666666; "G" = 0.02133733333333333; "H" = 0.064012
    "I" = 0.09402315767518445; "J" = 0.09402315767518446
    "K" = 3; "L" = 1
    "M" = 0.1112926666666667; "N" = 0.333878
    "O" = 0.01397697460904174; "P" = 0.01397697460904174
    "Q" = 0.002374688726222222; "R" = 0.021372198536
    "S" = 0.001314159287487981; "T" = 0.001314159287487981
  }
  6 = @{
    "A" = "MuSCs"; "B" = "Slit1"; "C" = "Gpc1"; "D" = "FAPs"
    "E" = 2; "F" = 0.6666666666666666; "G" = 0.02133733333333333; "H" = 0.064012
    "I" = 0.09402315767518445; "J" = 0.09402315767518446
    "K" = 3; "L" = 1
    "M" = 3.316850333333333; "N" = 9.950551000000001
    "O" = 0.4165551449121381; "P" = 0.4165551449121381
    "Q" = 0.07077274117911111; "R" = 0.6369546706120001
    "S" = 0.03916583007048327; "T" = 0.03916583007048328
  }
  7 = @{
    "A" = "MuSCs"; "B" = "Slit1"; "C" = "Gpc1"; "D" = "MuSCs"
    "E" = 2; "F" = 0.6666666666666666; "G" = 0.02133733333333333; "H" = 0.064012
    "I" = 0.09402315767518445; "J" = 0.09402315767518446
    "K" = 3; "L" = 1
    "M" = 4.534429; "N" = 13.603287
    "O" = 0.5694678804788202; "P" = 0.5694678804788201
    "Q" = 0.09675262304933334; "R" = 0.870773607444
    "S" = 0.05354316831721321; "T" = 0.0535431683172132
  }
}

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in $rows.Keys) {
  $rowData = $rows[$r]
  foreach ($c in $cols) {
    $ws.Range("$c$r").Value = $rowData[$c]
  }
}
